$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Q0)
$ws.Range("B3").Value = 0.2379616757713074
$ws.Range("C3").Value = 0.9704757615466515
$ws.Range("D3").Value = 1.752157255497987
$ws.Range("E3").Value = 1.323690770345546
$ws.Range("F3").Value = 1.331390667622807
$ws.Range("G3").Value = 23

# Row 4 (Q1)
$ws.Range("B4").Value = 0.2465549019837904
$ws.Range("C4").Value = 1.496450223635613
$ws.Range("D4").Value = 10.38184075314966
$ws.Range("E4").Value = 3.22208639753028
$ws.Range("F4").Value = 3.288241182326179
$ws.Range("G4").Value = 22
